$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns at B, shifting existing B:V data to K:AE.
$ws.Range("B1:J1").EntireColumn.Insert()

# New header dates (most-recent-first), left to right: Sep_08 .. Jun_16
$ws.Range("B1").Value2 = "Sep_08"
$ws.Range("C1").Value2 = "Aug_25"
$ws.Range("D1").Value2 = "Aug_04"
$ws.Range("E1").Value2 = "Jul_23"
$ws.Range("F1").Value2 = "Jul_17"
$ws.Range("G1").Value2 = "Jul_07"
$ws.Range("H1").Value2 = "Jun_30"
$ws.Range("I1").Value2 = "Jun_24"
$ws.Range("J1").Value2 = "Jun_16"

# Fill the new columns (B:J) with "UN" for every data row that has content,
# matching each row's existing extent.
for ($r = 2; $r -le 29; $r++) {
    $ws.Range("B$r`:J$r").Value2 = "UN"
}
for ($r = 30; $r -le 33; $r++) {
    $ws.Range("B$r`:J$r").Value2 = "UN"
}
